# Add season-record columns (Wins / Losses / Ties) to the sheet.
# The old scraper only pulled team statistics, not the season record, so
# this adds the missing columns and backfills every player row with the
# team's W-L-T record for the season.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# Columns AD (30), AE (31), AF (32) get new headers. Copy the formatting
# of the existing header cell (AC1) onto them first so they pick up the
# same bold/border/center style already used by the rest of row 1,
# then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Cells.Item(1, 30).Value = "Wins"
$ws.Cells.Item(1, 31).Value = "Losses"
$ws.Cells.Item(1, 32).Value = "Ties"

# --- Data rows ---------------------------------------------------------
# Every player (rows 2-52) shares the same team season record.
$wins = 74
$losses = 88
$ties = 0

for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}

Write-Output "season record columns added"
